# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# price table on Sheet1 to the new scraped values.
#
# D-column values that look numeric (e.g. "1.003") are written with a
# leading apostrophe so Excel keeps them as literal text, matching the
# original inline-string cells (e.g. "317.58") instead of silently
# coercing them into numbers and dropping trailing zeros / truncating
# the 3-decimal display.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.914.46'
$ws.Range("E2").Value = '  -3.65%  '
$ws.Range("D3").Value = '1.863.11'
$ws.Range("E3").Value = '  -2.81%  '
$ws.Range("D4").Value = '''1.003'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''317.49'
$ws.Range("E5").Value = '  -2.36%  '
$ws.Range("D6").Value = '''1.002'
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("D7").Value = '''0.4348'
$ws.Range("E7").Value = '  -5.41%  '
$ws.Range("D8").Value = '''0.3727'
$ws.Range("E8").Value = '  -2.49%  '
$ws.Range("D9").Value = '''0.07466'
$ws.Range("E9").Value = '  -3.16%  '
$ws.Range("D10").Value = '''0.9337'
$ws.Range("E10").Value = '  -4.85%  '
$ws.Range("D11").Value = '''21.31'
$ws.Range("D12").Value = '1.927.15'
$ws.Range("E12").Value = '  -0.05%  '
$ws.Range("D13").Value = '''6.731'
$ws.Range("E13").Value = '  -3.25%  '
$ws.Range("D14").Value = '''5.440'
$ws.Range("E14").Value = '  -4.29%  '
$ws.Range("D15").Value = '''0.06860'
$ws.Range("E15").Value = '  -1.79%  '
$ws.Range("E16").Value = '  +0.10%  '
$ws.Range("D17").Value = '''81.55'
$ws.Range("E17").Value = '  -3.10%  '
$ws.Range("E18").Value = '  -4.35%  '
$ws.Range("D19").Value = '''1.002'
$ws.Range("E19").Value = '  +0.15%  '
$ws.Range("D20").Value = '''15.86'
$ws.Range("E20").Value = '  -4.61%  '
$ws.Range("D21").Value = '27.906.39'
$ws.Range("E21").Value = '  -3.66%  '
$ws.Range("D23").Value = '''11.04'
$ws.Range("E23").Value = '  +0.91%  '
$ws.Range("D24").Value = '2.138.01'
$ws.Range("E24").Value = '  -0.70%  '
$ws.Range("E25").Value = '  -3.84%  '
$ws.Range("D26").Value = '''154.68'
$ws.Range("E26").Value = '  -2.34%  '
$ws.Range("D27").Value = '''18.44'
$ws.Range("E27").Value = '  -3.35%  '
$ws.Range("D28").Value = '''5.457'
$ws.Range("E28").Value = '  -4.03%  '
$ws.Range("D29").Value = '''113.35'
$ws.Range("E29").Value = '  -3.81%  '
$ws.Range("E30").Value = '  -7.53%  '
$ws.Range("D31").Value = '''0.09008'
$ws.Range("E31").Value = '  -3.11%  '
$ws.Range("D32").Value = '''0.8183'
$ws.Range("E32").Value = '  -5.41%  '
$ws.Range("D33").Value = '''4.813'
$ws.Range("E33").Value = '  -5.76%  '
$ws.Range("D34").Value = '''1.175'
$ws.Range("E34").Value = '  -6.17%  '
$ws.Range("D35").Value = '''2.971'
$ws.Range("E35").Value = '  -2.64%  '
$ws.Range("D36").Value = '''1.002'
$ws.Range("E36").Value = '  +0.17%  '
$ws.Range("E37").Value = '  -2.73%  '
$ws.Range("D38").Value = '''0.05490'
$ws.Range("E38").Value = '  -3.70%  '
$ws.Range("D39").Value = '''0.01973'
$ws.Range("E39").Value = '  -3.33%  '
$ws.Range("D40").Value = '''2.985'
$ws.Range("E40").Value = '  -2.88%  '
$ws.Range("D41").Value = '''0.5258'
$ws.Range("E41").Value = '  -4.52%  '
$ws.Range("D42").Value = '''7.033'
$ws.Range("E42").Value = '  -6.28%  '
$ws.Range("D43").Value = '''0.1702'
$ws.Range("E43").Value = '  -2.77%  '
$ws.Range("D44").Value = '''8.750'
$ws.Range("E44").Value = '  -6.69%  '
$ws.Range("D45").Value = '''0.06757'
$ws.Range("E45").Value = '  -2.06%  '
$ws.Range("D46").Value = '''0.4896'
$ws.Range("E46").Value = '  -5.50%  '
$ws.Range("D47").Value = '''10.65'
$ws.Range("E47").Value = '  -5.35%  '
$ws.Range("D48").Value = '''107.21'
$ws.Range("E48").Value = '  -2.90%  '
$ws.Range("D49").Value = '''1.675'
$ws.Range("E49").Value = '  -5.84%  '
$ws.Range("E50").Value = '  -0.01%  '
$ws.Range("D51").Value = '''1.887'
$ws.Range("E51").Value = '  -14.24%  '
